$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (ResNet-101, AdamW run 4): fill in the training/validation accuracy results
$ws.Range("J12").Value = 0.9162
$ws.Range("K12").Value = 0.8263

# Row 13 (ResNet-101, RMSprop run 44): fill in the training/validation accuracy results
$ws.Range("J13").Value = 0.9069
$ws.Range("K13").Value = 0.8279

# Row 20 (ResNet-50 run 44, epoch 20): fill in the training/validation accuracy results
$ws.Range("J20").Value = 0.9331
$ws.Range("K20").Value = 0.84

# Row 21 (ResNet-101, AdamW run 4, epoch 10): complete the whole row of data
$ws.Range("C21").Value = 0.0001
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 10
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "1e-3"
$ws.Range("G21").Value = 3
$ws.Range("H21").Value = 0.1
$ws.Range("J21").Value = 0.9414
$ws.Range("K21").Value = 0.8358
$ws.Range("L21").Value = "food_4_rn101.slurm / 24437"

# Row 22 (ResNet-101, RMSprop run 44, epoch 10): complete the whole row of data
$ws.Range("C22").Value = 0.0001
$ws.Range("D22").Value = 32
$ws.Range("E22").Value = 10
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "1e-3"
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 0.1
$ws.Range("J22").Value = 0.9289
$ws.Range("K22").Value = 0.8398
$ws.Range("L22").Value = "food_44_rn101.slurm / 24438"
